# Changes up to 22/11/23
# Update the "investments" sheet's Capex/FixOM figures for rows 13-16
# (M12-M15) and tidy up B14 so it matches the same number format/font
# as the rest of column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("investments")
$ws.Activate() | Out-Null

# Row 13 (M12): Capex 100000 -> 10000, FixOM 7000 -> 5000
$ws.Range("B13").Value = 10000
$ws.Range("C13").Value = 5000

# Row 14 (M13): Capex 66137 -> 15000, FixOM 576 -> 0
# B14 previously used a mismatched (11pt) font vs. the rest of the column
# (12pt) - bump it so it picks up the same style the other cells use.
$ws.Range("B14").Font.Size = 12
$ws.Range("B14").Value = 15000
$ws.Range("C14").Value = 0

# Row 15 (M14): Capex 100000 -> 10000, FixOM 7000 -> 5000
$ws.Range("B15").Value = 10000
$ws.Range("C15").Value = 5000

# Row 16 (M15): Capex 100000 -> 10000, FixOM 7000 -> 5000
$ws.Range("B16").Value = 10000
$ws.Range("C16").Value = 5000

# View state: zoomed out a bit and left the selection sitting on B15.
$excel.ActiveWindow.Zoom = 125
$ws.Range("B15").Select() | Out-Null
